$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.774.93'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '''1.894.91'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").Value = '''0.4941'
$ws.Range("E7").Value = '  +1.81%  '
$ws.Range("D8").Value = '''0.3795'
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").Value = '''0.07321'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = '''0.9109'
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("D11").Value = '''20.59'
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").Value = '''0.07617'
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("D13").Value = '''1.866.03'
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("D14").Value = '''5.464'
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").Value = '''91.11'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D18").Value = '''0.000008740'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '''27.806.28'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '''14.50'
$ws.Range("E21").Value = '  -3.49%  '
$ws.Range("D22").Value = '''5.118'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").Value = '''2.138.46'
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").Value = '''10.76'
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("D25").Value = '''154.00'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").Value = '''1.850'
$ws.Range("D27").Value = '''2.180'
$ws.Range("E27").Value = '  +3.62%  '
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("D29").Value = '''115.18'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").Value = '''4.882'
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("D31").Value = '''0.08937'
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").Value = '''3.250'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").Value = '''1.230'
$ws.Range("E33").Value = '  -1.88%  '
$ws.Range("D34").Value = '''0.7662'
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").Value = '''4.642'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '''0.02046'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").Value = '''2.552'
$ws.Range("E37").Value = '  -7.93%  '
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = '''0.5500'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("D40").Value = '''0.05285'
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").Value = '''2.988'
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("D42").Value = '''6.901'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").Value = '''8.551'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.1520'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '''112.51'
$ws.Range("E45").Value = '  +4.77%  '
$ws.Range("D46").Value = '''10.62'
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").Value = '''0.4788'
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("D48").Value = '''0.9996'
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").Value = '''1.633'
$ws.Range("D50").Value = '''67.42'
$ws.Range("E50").Value = '  -3.12%  '
$ws.Range("D51").Value = '''0.06061'
$ws.Range("E51").Value = '  -1.52%  '
